$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

$ws.Range("B7").Value = "yes"
$ws.Range("C7").ClearContents()
